$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: mark Invalid (G) and Absent (H)
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-11: mark Absent (H) only
foreach ($r in 4..11) {
    $ws.Range("H$r").Value = 1
}

# Row 12: mark Total Attendance Count (D) and Real (E)
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Rows 13-18: mark Absent (H) only
foreach ($r in 13..18) {
    $ws.Range("H$r").Value = 1
}
